# Append a new row (row 60) of sensor data to every worksheet in the
# workbook (ROW35-FE-LIFTER, ROW35-MID-LIFTER, ROW02-FE-LIFTER,
# ROW02-MID-LIFTER). Each sheet currently ends at row 59; this adds one
# more reading after it, growing the used range to A1:I60.

$wb = $excel.ActiveWorkbook

# Per-sheet values for the new row, in sheet (tab) order.
$rows = @(
    @{ A = "2025-03-06 19:42:06"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"; D = "0x01,0x90,"; E = "0x d";  F = 400; G = "568631262647113770877196"; H = 400; I = 13 },
    @{ A = "2025-03-06 19:29:35"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; D = "0x01,0x90,"; E = "0x e";  F = 400; G = "568631262647113770942732"; H = 400; I = 14 },
    @{ A = "2025-03-06 19:51:45"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"; D = "0x01,0x90,"; E = "0xff";   F = 400; G = "568631262647113769959692"; H = 400; I = 255 },
    @{ A = "2025-03-06 19:41:15"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x01,0x90,"; E = "0x 3";  F = 400; G = "568631262647113769959692"; H = 400; I = 3 }
)

for ($i = 0; $i -lt $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $data = $rows[$i]
    $newRow = $ws.UsedRange.Rows.Count + 1

    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F

    # Column G holds a 24-digit integer that exceeds double precision, so it
    # must stay a text value instead of being coerced into scientific
    # notation. Force text format for the assignment, then restore the
    # default (un-styled) cell style so the saved cell matches the rest of
    # the column, which carries no explicit style.
    $gCell = $ws.Cells.Item($newRow, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $data.G
    $gCell.Style = "Normal"

    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
